# Fill in May 2023 data on the "May" sheet (sheet6) and let the
# "Yearly totals" sheet's existing SUM formulas across January..December
# recompute automatically once May's raw numbers are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May")

$rows = @(
    @{Row=2; B=1613; C=1219; D=394; E='We borrowerd more than we lent'; F=$null; G='1.32 : 1'},
    @{Row=3; B=476; C=463; D=13; E='We borrowerd more than we lent'; F=$null; G='1.03 : 1'},
    @{Row=4; B=1279; C=1322; D=-43; E=$null; F='We lent more than we borrowed'; G='0.97 : 1'},
    @{Row=5; B=132; C=119; D=13; E='We borrowerd more than we lent'; F=$null; G='1.11 : 1'},
    @{Row=6; B=1218; C=1544; D=-326; E=$null; F='We lent more than we borrowed'; G='0.79 : 1'},
    @{Row=7; B=195; C=188; D=7; E='We borrowerd more than we lent'; F=$null; G='1.04 : 1'},
    @{Row=8; B=122; C=190; D=-68; E=$null; F='We lent more than we borrowed'; G='0.64 : 1'},
    @{Row=9; B=41; C=71; D=-30; E=$null; F='We lent more than we borrowed'; G='0.58 : 1'},
    @{Row=10; B=5; C=31; D=-26; E=$null; F='We lent more than we borrowed'; G='0.16 : 1'},
    @{Row=11; B=0; C=0; D=0; E=$null; F=$null; G=$null},
    @{Row=12; B=8; C=10; D=-2; E=$null; F='We lent more than we borrowed'; G='0.80 : 1'},
    @{Row=13; B=171; C=73; D=98; E='We borrowerd more than we lent'; F=$null; G='2.34 : 1'},
    @{Row=14; B=111; C=273; D=-162; E=$null; F='We lent more than we borrowed'; G='0.41 : 1'},
    @{Row=15; B=77; C=109; D=-32; E=$null; F='We lent more than we borrowed'; G='0.71 : 1'},
    @{Row=16; B=61; C=165; D=-104; E=$null; F='We lent more than we borrowed'; G='0.37 : 1'},
    @{Row=17; B=555; C=454; D=101; E='We borrowerd more than we lent'; F=$null; G='1.22 : 1'},
    @{Row=18; B=89; C=89; D=0; E=$null; F=$null; G='1.00 : 1'},
    @{Row=19; B=707; C=419; D=288; E='We borrowerd more than we lent'; F=$null; G='1.69 : 1'},
    @{Row=20; B=3; C=74; D=-71; E=$null; F='We lent more than we borrowed'; G='0.04 : 1'},
    @{Row=21; B=496; C=420; D=76; E='We borrowerd more than we lent'; F=$null; G='1.18 : 1'},
    @{Row=22; B=32; C=43; D=-11; E=$null; F='We lent more than we borrowed'; G='0.74 : 1'},
    @{Row=23; B=763; C=519; D=244; E='We borrowerd more than we lent'; F=$null; G='1.47 : 1'},
    @{Row=24; B=1907; C=1256; D=651; E='We borrowerd more than we lent'; F=$null; G='1.52 : 1'},
    @{Row=25; B=176; C=327; D=-151; E=$null; F='We lent more than we borrowed'; G='0.54 : 1'},
    @{Row=26; B=0; C=0; D=0; E=$null; F=$null; G=$null},
    @{Row=27; B=209; C=211; D=-2; E=$null; F='We lent more than we borrowed'; G='0.99 : 1'},
    @{Row=28; B=53; C=62; D=-9; E=$null; F='We lent more than we borrowed'; G='0.85 : 1'},
    @{Row=29; B=449; C=445; D=4; E='We borrowerd more than we lent'; F=$null; G='1.01 : 1'},
    @{Row=30; B=18; C=50; D=-32; E=$null; F='We lent more than we borrowed'; G='0.36 : 1'},
    @{Row=31; B=61; C=299; D=-238; E=$null; F='We lent more than we borrowed'; G='0.20 : 1'},
    @{Row=32; B=440; C=572; D=-132; E=$null; F='We lent more than we borrowed'; G='0.77 : 1'},
    @{Row=33; B=363; C=497; D=-134; E=$null; F='We lent more than we borrowed'; G='0.73 : 1'},
    @{Row=34; B=171; C=102; D=69; E='We borrowerd more than we lent'; F=$null; G='1.68 : 1'},
    @{Row=35; B=880; C=1019; D=-139; E=$null; F='We lent more than we borrowed'; G='0.86 : 1'},
    @{Row=36; B=241; C=421; D=-180; E=$null; F='We lent more than we borrowed'; G='0.57 : 1'},
    @{Row=37; B=467; C=359; D=108; E='We borrowerd more than we lent'; F=$null; G='1.30 : 1'},
    @{Row=38; B=11; C=166; D=-155; E=$null; F='We lent more than we borrowed'; G='0.07 : 1'},
    @{Row=39; B=1; C=48; D=-47; E=$null; F='We lent more than we borrowed'; G='0.02 : 1'},
    @{Row=40; B=15; C=73; D=-58; E=$null; F='We lent more than we borrowed'; G='0.21 : 1'},
    @{Row=41; B=3; C=20; D=-17; E=$null; F='We lent more than we borrowed'; G='0.15 : 1'},
    @{Row=42; B=0; C=1; D=-1; E=$null; F='We lent more than we borrowed'; G='0.00 : 1'},
    @{Row=44; B=51; C=68; D=-17; E=$null; F='We lent more than we borrowed'; G='0.75 : 1'},
    @{Row=45; B=77; C=153; D=-76; E=$null; F='We lent more than we borrowed'; G='0.50 : 1'},
    @{Row=46; B=401; C=562; D=-161; E=$null; F='We lent more than we borrowed'; G='0.71 : 1'},
    @{Row=47; B=959; C=538; D=421; E='We borrowerd more than we lent'; F=$null; G='1.78 : 1'},
    @{Row=48; B=231; C=676; D=-445; E=$null; F='We lent more than we borrowed'; G='0.34 : 1'},
    @{Row=49; B=311; C=178; D=133; E='We borrowerd more than we lent'; F=$null; G='1.75 : 1'},
    @{Row=50; B=815; C=533; D=282; E='We borrowerd more than we lent'; F=$null; G='1.53 : 1'},
    @{Row=51; B=293; C=66; D=227; E='We borrowerd more than we lent'; F=$null; G='4.44 : 1'},
    @{Row=52; B=475; C=507; D=-32; E=$null; F='We lent more than we borrowed'; G='0.94 : 1'},
    @{Row=53; B=136; C=238; D=-102; E=$null; F='We lent more than we borrowed'; G='0.57 : 1'},
    @{Row=54; B=22; C=243; D=-221; E=$null; F='We lent more than we borrowed'; G='0.09 : 1'},
    @{Row=55; B=278; C=183; D=95; E='We borrowerd more than we lent'; F=$null; G='1.52 : 1'}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B   # B: Other Next materials checked out at our library
    $ws.Cells.Item($row, 3).Value = $r.C   # C: Our materials checked out at other Next libraries
    $ws.Cells.Item($row, 4).Value = $r.D   # D: Net (B - C)

    if ($r.E -ne $null) {
        $ws.Cells.Item($row, 5).Value = $r.E
    }
    if ($r.F -ne $null) {
        $ws.Cells.Item($row, 6).Value = $r.F
    }
    if ($r.G -ne $null) {
        $ws.Cells.Item($row, 7).Value = $r.G
    }
}
